$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits right after the
#    "...guardar los cambios en el codigo fuente)" run.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Type "Modules:" into what is currently the last (empty) paragraph of the
#    document.
# ---------------------------------------------------------------------------
$modulesPara = $d.Paragraphs.Last
$modulesPara.Range.InsertAfter("Modules:")
$modulesPara.Range.LanguageID = "es-419"

# ---------------------------------------------------------------------------
# 3. Append a brand-new paragraph after it describing AppModule. This mirrors
#    what Word itself would have produced while the user was mid-edit
#    (word split into multiple runs because of the underline formatting and
#    the then-current "_GoBack" cursor position, plus proofing marks around
#    the two words flagged by the spell checker).
# ---------------------------------------------------------------------------
$insertAt = $d.Range($d.Content.End, $d.Content.End)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="es-419"/>
              </w:rPr>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
                <w:lang w:val="es-419"/>
              </w:rPr>
              <w:t>App</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
                <w:lang w:val="es-419"/>
              </w:rPr>
              <w:t>M</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:u w:val="single"/>
                <w:lang w:val="es-419"/>
              </w:rPr>
              <w:t>odule</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-419"/>
              </w:rPr>
              <w:t xml:space="preserve"> es la ra&#237;z de la </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-419"/>
              </w:rPr>
              <w:t>aplicacion</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$insertAt.InsertXML($xml)
